$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.19105733333333
$ws.Range("H2").Value = 48.573172
$ws.Range("I2").Value = 0.0401918797050022
$ws.Range("J2").Value = 0.0401918797050022
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 11.25749966666667
$ws.Range("N2").Value = 33.772499
$ws.Range("O2").Value = 0.6929800609896341
$ws.Range("P2").Value = 0.6929800609896341
$ws.Range("Q2").Value = 182.2708225329809
$ws.Range("R2").Value = 1640.437402796828
$ws.Range("S2").Value = 0.02785217124926046
$ws.Range("T2").Value = 0.02785217124926046

$ws.Range("G3").Value = 16.19105733333333
$ws.Range("H3").Value = 48.573172
$ws.Range("I3").Value = 0.0401918797050022
$ws.Range("J3").Value = 0.0401918797050022
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9898276666666668
$ws.Range("N3").Value = 2.969483
$ws.Range("O3").Value = 0.06093101107050686
$ws.Range("P3").Value = 0.06093101107050686
$ws.Range("Q3").Value = 16.02635650111956
$ws.Range("R3").Value = 144.237208510076
$ws.Range("S3").Value = 0.002448931867249969
$ws.Range("T3").Value = 0.002448931867249969

$ws.Range("G4").Value = 16.19105733333333
$ws.Range("H4").Value = 48.573172
$ws.Range("I4").Value = 0.0401918797050022
$ws.Range("J4").Value = 0.0401918797050022
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.821582
$ws.Range("N4").Value = 11.464746
$ws.Range("O4").Value = 0.2352458543950409
$ws.Range("P4").Value = 0.2352458543950409
$ws.Range("Q4").Value = 61.87545326603466
$ws.Range("R4").Value = 556.879079394312
$ws.Range("S4").Value = 0.009454973080945946
$ws.Range("T4").Value = 0.009454973080945946

$ws.Range("G5").Value = 16.19105733333333
$ws.Range("H5").Value = 48.573172
$ws.Range("I5").Value = 0.0401918797050022
$ws.Range("J5").Value = 0.0401918797050022
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1761463333333333
$ws.Range("N5").Value = 0.528439
$ws.Range("O5").Value = 0.01084307354481826
$ws.Range("P5").Value = 0.01084307354481827
$ws.Range("Q5").Value = 2.851995382056444
$ws.Range("R5").Value = 25.667958438508
$ws.Range("S5").Value = 0.0004358035075458274
$ws.Range("T5").Value = 0.0004358035075458274

$ws.Range("G6").Value = 255.6993613333333
$ws.Range("H6").Value = 767.098084
$ws.Range("I6").Value = 0.6347354443738135
$ws.Range("J6").Value = 0.6347354443738134
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 11.25749966666667
$ws.Range("N6").Value = 33.772499
$ws.Range("O6").Value = 0.6929800609896341
$ws.Range("P6").Value = 0.6929800609896341
$ws.Range("Q6").Value = 2878.535474976879
$ws.Range("R6").Value = 25906.81927479191
$ws.Range("S6").Value = 0.4398590069544477
$ws.Range("T6").Value = 0.4398590069544476

$ws.Range("G7").Value = 255.6993613333333
$ws.Range("H7").Value = 767.098084
$ws.Range("I7").Value = 0.6347354443738135
$ws.Range("J7").Value = 0.6347354443738134
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9898276666666668
$ws.Range("N7").Value = 2.969483
$ws.Range("O7").Value = 0.06093101107050686
$ws.Range("P7").Value = 0.06093101107050686
$ws.Range("Q7").Value = 253.0983021967302
$ws.Range("R7").Value = 2277.884719770572
$ws.Range("S7").Value = 0.03867507238798392
$ws.Range("T7").Value = 0.03867507238798391

$ws.Range("G8").Value = 255.6993613333333
$ws.Range("H8").Value = 767.098084
$ws.Range("I8").Value = 0.6347354443738135
$ws.Range("J8").Value = 0.6347354443738134
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.821582
$ws.Range("N8").Value = 11.464746
$ws.Range("O8").Value = 0.2352458543950409
$ws.Range("P8").Value = 0.2352458543950409
$ws.Range("Q8").Value = 977.1760766829625
$ws.Range("R8").Value = 8794.584690146663
$ws.Range("S8").Value = 0.1493188819265337
$ws.Range("T8").Value = 0.1493188819265337

$ws.Range("G9").Value = 255.6993613333333
$ws.Range("H9").Value = 767.098084
$ws.Range("I9").Value = 0.6347354443738135
$ws.Range("J9").Value = 0.6347354443738134
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1761463333333333
$ws.Range("N9").Value = 0.528439
$ws.Range("O9").Value = 0.01084307354481826
$ws.Range("P9").Value = 0.01084307354481827
$ws.Range("Q9").Value = 45.04050493454177
$ws.Range("R9").Value = 405.364544410876
$ws.Range("S9").Value = 0.00688248310484816
$ws.Range("T9").Value = 0.006882483104848161

$ws.Range("G10").Value = 57.51647566666667
$ws.Range("H10").Value = 172.549427
$ws.Range("I10").Value = 0.142776053685583
$ws.Range("J10").Value = 0.142776053685583
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 11.25749966666667
$ws.Range("N10").Value = 33.772499
$ws.Range("O10").Value = 0.6929800609896341
$ws.Range("P10").Value = 0.6929800609896341
$ws.Range("Q10").Value = 647.4917056453414
$ws.Range("R10").Value = 5827.425350808073
$ws.Range("S10").Value = 0.0989409583908946
$ws.Range("T10").Value = 0.09894095839089459

$ws.Range("G11").Value = 57.51647566666667
$ws.Range("H11").Value = 172.549427
$ws.Range("I11").Value = 0.142776053685583
$ws.Range("J11").Value = 0.142776053685583
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.9898276666666668
$ws.Range("N11").Value = 2.969483
$ws.Range("O11").Value = 0.06093101107050686
$ws.Range("P11").Value = 0.06093101107050686
$ws.Range("Q11").Value = 56.93139890402679
$ws.Range("R11").Value = 512.3825901362411
$ws.Range("S11").Value = 0.008699489307719542
$ws.Range("T11").Value = 0.008699489307719542

$ws.Range("G12").Value = 57.51647566666667
$ws.Range("H12").Value = 172.549427
$ws.Range("I12").Value = 0.142776053685583
$ws.Range("J12").Value = 0.142776053685583
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.821582
$ws.Range("N12").Value = 11.464746
$ws.Range("O12").Value = 0.2352458543950409
$ws.Range("P12").Value = 0.2352458543950409
$ws.Range("Q12").Value = 219.8039281111713
$ws.Range("R12").Value = 1978.235353000542
$ws.Range("S12").Value = 0.03358747473641721
$ws.Range("T12").Value = 0.0335874747364172

$ws.Range("G13").Value = 57.51647566666667
$ws.Range("H13").Value = 172.549427
$ws.Range("I13").Value = 0.142776053685583
$ws.Range("J13").Value = 0.142776053685583
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1761463333333333
$ws.Range("N13").Value = 0.528439
$ws.Range("O13").Value = 0.01084307354481826
$ws.Range("P13").Value = 0.01084307354481827
$ws.Range("Q13").Value = 10.13131629493922
$ws.Range("R13").Value = 91.181846654453
$ws.Range("S13").Value = 0.001548131250551697
$ws.Range("T13").Value = 0.001548131250551698

$ws.Range("G14").Value = 73.43709933333334
$ws.Range("H14").Value = 220.311298
$ws.Range("I14").Value = 0.1822966222356014
$ws.Range("J14").Value = 0.1822966222356014
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 11.25749966666667
$ws.Range("N14").Value = 33.772499
$ws.Range("O14").Value = 0.6929800609896341
$ws.Range("P14").Value = 0.6929800609896341
$ws.Range("Q14").Value = 826.7181212659668
$ws.Range("R14").Value = 7440.463091393701
$ws.Range("S14").Value = 0.1263279243950313
$ws.Range("T14").Value = 0.1263279243950313

$ws.Range("G15").Value = 73.43709933333334
$ws.Range("H15").Value = 220.311298
$ws.Range("I15").Value = 0.1822966222356014
$ws.Range("J15").Value = 0.1822966222356014
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.9898276666666668
$ws.Range("N15").Value = 2.969483
$ws.Range("O15").Value = 0.06093101107050686
$ws.Range("P15").Value = 0.06093101107050686
$ws.Range("Q15").Value = 72.69007267988157
$ws.Range("R15").Value = 654.210654118934
$ws.Range("S15").Value = 0.01110751750755344
$ws.Range("T15").Value = 0.01110751750755343

$ws.Range("G16").Value = 73.43709933333334
$ws.Range("H16").Value = 220.311298
$ws.Range("I16").Value = 0.1822966222356014
$ws.Range("J16").Value = 0.1822966222356014
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.821582
$ws.Range("N16").Value = 11.464746
$ws.Range("O16").Value = 0.2352458543950409
$ws.Range("P16").Value = 0.2352458543950409
$ws.Range("Q16").Value = 280.6458969444786
$ws.Range("R16").Value = 2525.813072500308
$ws.Range("S16").Value = 0.04288452465114406
$ws.Range("T16").Value = 0.04288452465114406

$ws.Range("G17").Value = 73.43709933333334
$ws.Range("H17").Value = 220.311298
$ws.Range("I17").Value = 0.1822966222356014
$ws.Range("J17").Value = 0.1822966222356014
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1761463333333333
$ws.Range("N17").Value = 0.528439
$ws.Range("O17").Value = 0.01084307354481826
$ws.Range("P17").Value = 0.01084307354481827
$ws.Range("Q17").Value = 12.93567577820244
$ws.Range("R17").Value = 116.421082003822
$ws.Range("S17").Value = 0.001976655681872578
$ws.Range("T17").Value = 0.001976655681872578

